$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F ("dSF")
$updates = @{
    4  = -6
    9  = -5
    12 = -10
    14 = -5
    15 = -3
    17 = 9
    18 = -1
    21 = 4
    27 = -1
    28 = 0
    46 = -7
    49 = 1
    53 = 13
    57 = 0
    59 = 0
    63 = 7
    64 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
